# This script implements the weekly update to the "Espinaca" price table.
# A new weekly record is inserted at row 40 (pushing all subsequent rows
# down by one), and the record that used to be the very last data row
# (row 164) is appended as the new last row (165).
#
# Only columns D (Fecha), J (Volumen), K (Precio minimo), M (Precio
# promedio ponderado) and P (Precio $/Kg) vary from row to row - every
# other column holds a value that is constant for the whole table, so
# those columns do not need to be touched for the newly appended row;
# we simply copy them down from the row above.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 40
$lastDataRow = 164
$newLastRow = 165

# --- Step 1: snapshot the original values of the variable columns for
#     every row that is going to move, BEFORE any writes happen.
$origD = @{}
$origJ = @{}
$origK = @{}
$origM = @{}
$origP = @{}

for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $origD[$r] = $ws.Cells.Item($r, 4).Value2()
    $origJ[$r] = $ws.Cells.Item($r, 10).Value2()
    $origK[$r] = $ws.Cells.Item($r, 11).Value2()
    $origM[$r] = $ws.Cells.Item($r, 13).Value2()
    $origP[$r] = $ws.Cells.Item($r, 16).Value2()
}

# --- Step 2: append the new last row (165), duplicating every constant
#     column from the row above it, then set the variable columns to
#     what used to be row 164's values (the record that is being pushed
#     off the end of the shifted block). Column D carries a date number
#     format, so its format is copied explicitly too.
for ($c = 1; $c -le 18; $c++) {
    $ws.Cells.Item($newLastRow, $c).Value = $ws.Cells.Item($lastDataRow, $c).Value2()
}
$ws.Cells.Item($newLastRow, 4).NumberFormat = $ws.Cells.Item($lastDataRow, 4).NumberFormat

$ws.Cells.Item($newLastRow, 4).Value2 = $origD[$lastDataRow]
$ws.Cells.Item($newLastRow, 10).Value2 = $origJ[$lastDataRow]
$ws.Cells.Item($newLastRow, 11).Value2 = $origK[$lastDataRow]
$ws.Cells.Item($newLastRow, 13).Value2 = $origM[$lastDataRow]
$ws.Cells.Item($newLastRow, 16).Value2 = $origP[$lastDataRow]

# --- Step 3: shift every row from 164 down to 41 so it takes on the
#     values that used to belong to the row directly above it.
for ($r = $lastDataRow; $r -ge ($firstDataRow + 1); $r--) {
    $ws.Cells.Item($r, 4).Value2 = $origD[$r - 1]
    $ws.Cells.Item($r, 10).Value2 = $origJ[$r - 1]
    $ws.Cells.Item($r, 11).Value2 = $origK[$r - 1]
    $ws.Cells.Item($r, 13).Value2 = $origM[$r - 1]
    $ws.Cells.Item($r, 16).Value2 = $origP[$r - 1]
}

# --- Step 4: write the brand-new weekly record into row 40.
$ws.Cells.Item($firstDataRow, 4).Value2 = 44487
$ws.Cells.Item($firstDataRow, 10).Value2 = 3300
$ws.Cells.Item($firstDataRow, 11).Value2 = 400
$ws.Cells.Item($firstDataRow, 13).Value2 = 450
$ws.Cells.Item($firstDataRow, 16).Value2 = 900
